# 15HP01_index.xlsx -- header normalisation
#
# The sheet's header row (A1:H1) is re-labelled to lower-case column
# names (Experiment -> experiment, Plant_Code -> plant_code, etc.) and
# the header row is selected (A1:H1) instead of just the first cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "experiment"
$ws.Range("B1").Value = "plant_code"
$ws.Range("C1").Value = "genotype"
$ws.Range("D1").Value = "position"
$ws.Range("E1").Value = "line"
$ws.Range("F1").Value = "column"
$ws.Range("G1").Value = "repeat"
$ws.Range("H1").Value = "treatment"

$ws.Range("A1:H1").Select() | Out-Null

Write-Output "Header row re-labelled (lower-case) and A1:H1 selected."
